# Apply updated K (strikeout) values to column G, rows 2-38 of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 2
    3 = 3
    4 = 2
    5 = 2
    6 = 3
    7 = 3
    8 = 0
    9 = 2
    10 = 2
    11 = 4
    12 = 0
    13 = 2
    14 = 1
    15 = 2
    16 = 0
    17 = 2
    18 = 2
    19 = 6
    20 = 1
    21 = 6
    22 = 3
    23 = 3
    24 = 2
    25 = 0
    26 = 3
    27 = 0
    28 = 1
    29 = 3
    30 = 1
    31 = 1
    32 = 4
    33 = 1
    34 = 2
    35 = 5
    36 = 4
    37 = 3
    38 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G" + $row).Value = $gValues[$row]
}

